# Planification.xlsx - "Etat de l'art II - comparaison" commit
# Adds the 4 tasks logged on 2020-04-21 to the "SPRINT 3" sheet, extends the
# SUM() range that totals the sprint's hours, and leaves the selection on
# the first empty row below the new entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SPRINT 3")

# Row 12: Test Copyleaks API (1h)
$ws.Range("B12").Value = "2020-04-21"
$ws.Range("C12").Value = "Test Copyleaks API"
$ws.Range("D12").Value = 1

# Row 13: Test connection FTP server dev (0.5h)
$ws.Range("B13").Value = "2020-04-21"
$ws.Range("C13").Value = "Test connection FTP server dev"
$ws.Range("D13").Value = 0.5

# Row 14: Comparaison des outils, calcul prix (1h)
$ws.Range("B14").Value = "2020-04-21"
$ws.Range("C14").Value = "Comparaison des outils, calcul prix"
$ws.Range("D14").Value = 1

# Row 15: Planification (0.5h)
$ws.Range("B15").Value = "2020-04-21"
$ws.Range("C15").Value = "Planification"
$ws.Range("D15").Value = 0.5

# Extend the hours total to cover the newly-added rows
$ws.Range("D17").Formula = "=SUM(D3:D15)"

# Leave the cursor on the next empty row, like the author did
$ws.Range("B16").Select() | Out-Null
